$d = $word.ActiveDocument

# --- Step 1: "Project 3" -> "Project 5" ---
$p4 = $d.Paragraphs.Item(4)
$r4 = $p4.Range
$r4.MoveEnd(1, -1)
[void]$r4.Find.Execute('3', $false, $false, $false, $false, $false, $true, 1, $false, '5', 2)

# --- Step 2: Delete paragraph 8 (old "extended ... 22 rows ... 76 columns" paragraph) entirely ---
$p8 = $d.Paragraphs.Item(8)
$p8.Range.Delete()

# --- Step 3: Delete paragraph 7 (old "recursion" paragraph) entirely ---
$p7 = $d.Paragraphs.Item(7)
$p7.Range.Delete()

# --- Step 4: Replace paragraph 6 (ifstream paragraph) content with the new "List class" text ---
$p6 = $d.Paragraphs.Item(6)
$r6 = $p6.Range
$r6.MoveEnd(1, -1)
$r6.Text = 'The program creates one List class object, two integer variables and one char variable. It shows the user the menu with options to choose from and asks the user to choose an option from the list. It read the user''s input and apply the action the user has chosen on the list, the program after that shows the menu options again and asks the user to choose an action. It repeats this process until the user decides to quit the program. The program also checks if the user has entered a wrong menu option and asks the user to renter a correct menu option every time the user enters a wrong menu option.'

# Wrap the whole paragraph (content + paragraph mark) with bookmark _Hlk506402853
$p6b = $d.Paragraphs.Item(6)
$bmRange = $d.Range($p6b.Range.Start, $p6b.Range.End)
$d.Bookmarks.Add('_Hlk506402853', $bmRange)

# --- Step 5: Insert new paragraph after paragraph 6 for "Ordered lists..." ---
$p6c = $d.Paragraphs.Item(6)
$p6c.Range.InsertParagraphAfter()
$pOrdered = $d.Paragraphs.Item(7)
$rOrdered = $pOrdered.Range
$rOrdered.MoveEnd(1, -1)
$rOrdered.Text = 'Ordered lists are lists that holds items all arranged according to the value of each item after the other. Linked List data structure is useful to implement ordered lists. This program uses an ordered list ADT in a class with a linked list. The program uses this class to get offer a broad range of actions to perform on the list.'

# --- Step 6: Insert new paragraph after paragraph 7 for the "extended...Stacks and Queues" text ---
$p7c = $d.Paragraphs.Item(7)
$p7c.Range.InsertParagraphAfter()
$pExt = $d.Paragraphs.Item(8)
$rExt = $pExt.Range
$rExt.MoveEnd(1, -1)
$rExt.Text = 'The program can be extended in many ways, we can add functions that allow merging two lists, add more functions to allow the user to add different items in the list but they won’t be ordered, we can also use different data structures as Stacks and Queues to implement different versions of the program. '

# Insert collapsed _GoBack bookmark right before the final ". " run
$pExt2 = $d.Paragraphs.Item(8)
$gbPos = $pExt2.Range.End - 1 - 2
$gbRange = $d.Range($gbPos, $gbPos)
$d.Bookmarks.Add('_GoBack', $gbRange)
